# Weekly update: insert a new latest-week row at the top of the data
# (row 145), pushing the existing rows 145:166 down to 146:167.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 145:166 down by one row, creating a fresh blank row 145.
$ws.Rows("145:145").Insert()

# Populate the newly inserted row 145 with this week's record.
$ws.Cells.Item(145, 1).Value = 7
$ws.Cells.Item(145, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(145, 3).Value = "Ñuble"
$ws.Cells.Item(145, 4).Value = 44474
$ws.Cells.Item(145, 5).Value = 16
$ws.Cells.Item(145, 6).Value = 100112008
$ws.Cells.Item(145, 7).Value = "Coliflor"
$ws.Cells.Item(145, 8).Value = "Sin especificar"
$ws.Cells.Item(145, 9).Value = "Primera"
$ws.Cells.Item(145, 10).Value = 120
$ws.Cells.Item(145, 11).Value = 650
$ws.Cells.Item(145, 12).Value = 700
$ws.Cells.Item(145, 13).Value = 675
$ws.Cells.Item(145, 14).Value = "`$/unidad"
$ws.Cells.Item(145, 15).Value = "Región del Maule"
$ws.Cells.Item(145, 16).Value = 675
$ws.Cells.Item(145, 17).Value = 1
$ws.Cells.Item(145, 18).Value = "Hortaliza"
